# Actualización automática 2025-09-15 09:54:22
#
# Inserts a new client row ("GONZALEZ CARDENAS ERNESTO PAOLO") in alphabetical
# position (row 15) on both the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets.
# This shifts every subsequent client row down by one and pushes the trailing
# totals/summary row from 31 to 32. The summary row on "VENTAS POR GRUPO"
# contains "N de 29" style labels that must be updated to "N de 30" to reflect
# the new total client count.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "VENTAS POR GRUPO" (columns A..R) ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Push row 15 and everything below it (including the trailing totals row)
# down by one row, duplicating formatting as Excel normally does on insert.
$ws1.Rows.Item(15).Insert()

# Populate the freshly inserted, now-blank row 15 for the new client.
$ws1.Range("A15").Value = "LOZANO MOLINA TITO"
$ws1.Range("B15").Value = "GONZALEZ CARDENAS ERNESTO PAOLO"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(15, $col).Value = 0
}

# --- Sheet 2: "VENTA MENSUAL" (columns A..G) ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(15).Insert()

$ws2.Range("A15").Value = "LOZANO MOLINA TITO"
$ws2.Range("B15").Value = "GONZALEZ CARDENAS ERNESTO PAOLO"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(15, $col).Value = 0
}

# Update the "X de 29" -> "X de 30" labels in the summary row of
# "VENTAS POR GRUPO", which is now row 32 after the insert above.
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(32, $col)
    $text = $cell.Value()
    if ($text -ne $null) {
        $cell.Value = $text -replace "de 29", "de 30"
    }
}
